$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "character_image_path" column (J) is being dropped from the
# CharacterData sheet - clear out the header + all per-character values.
$rng = $ws.Range("J1:J8")
$rng.ClearContents()
$rng.NumberFormat = "General"

# Move the active selection to where the author left it (J10).
$ws.Range("J10").Select()
